$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) from 45202 to 45203 for all existing data rows (2-270)
$ws.Range("C2:C270").Value = 45203

# Row 270 picks up an explicit row height once the new row is appended below it
$ws.Rows.Item(270).RowHeight = 15

# Add the new row 271 (A 47152-2023)
$row = 271
$ws.Cells.Item($row, 1).Value = "A 47152-2023"
$ws.Cells.Item($row, 2).Value = 45201
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row, 3).Value = 45203
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item($row, 5).Value = "GRUMS"
$ws.Cells.Item($row, 7).Value = 0.5
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 18).WrapText = $true
